# "Start of data wrangling" - clean up the Asthma by state sheet:
#  - Rename headers: "State or Territory" -> "State",
#    "Percent With Current Asthma (SE)" -> "Percent With Current Asthma"
#  - Convert the Column C "xx.x (yy.yy)" strings into plain numeric percentages
#    (drop the standard-error part in parentheses), leaving footnoted cells
#    (e.g. Florida's missing-data marker) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "State"
$ws.Range("C1").Value = "Percent With Current Asthma"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($val -ne $null -and $val.GetType().Name -eq "String") {
        $parts = $val.Split(" ")
        $numPart = $parts[0]
        $num = $null
        try {
            $candidate = $numPart + 0.0
            if ($candidate.GetType().Name -ne "String") {
                $num = $candidate
            }
        } catch {
            $num = $null
        }
        if ($num -ne $null) {
            $cell.Value = $num
        }
    }
}

$ws.Range("E15").Select()
